$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.196.78"
$ws.Range("E2").Value = "  +2.27%  "

# Row 3
$ws.Range("D3").Value = "2.529.89"
$ws.Range("E3").Value = "  +1.62%  "

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'596.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "

# Row 6
$ws.Range("D6").Value = "'177.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "

# Row 7
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").Value = "'0.520"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.03%  "

# Row 9
$ws.Range("D9").Value = "2.528.91"
$ws.Range("E9").Value = "  +1.53%  "

# Row 10
$ws.Range("D10").Value = "'0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.68%  "

# Row 11
$ws.Range("E11").Value = "  -1.04%  "

# Row 12
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'4.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.29%  "

# Row 13
$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D13").Value = "'0.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.74%  "

# Row 14
$ws.Range("D14").Value = "'26.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.11%  "

# Row 15
$ws.Range("D15").Value = "2.951.36"

# Row 16
$ws.Range("D16").Value = "68.906.35"
$ws.Range("E16").Value = "  +2.00%  "

# Row 17
$ws.Range("E17").Value = "  +1.08%  "

# Row 18
$ws.Range("D18").Value = "2.523.94"
$ws.Range("E18").Value = "  -1.19%  "

# Row 20
$ws.Range("D20").Value = "'361.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.94%  "

# Row 21
$ws.Range("D21").Value = "'7.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.36%  "

# Row 22
$ws.Range("D22").Value = "'4.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.83%  "

# Row 23
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("D24").Value = "'70.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.63%  "

# Row 25
$ws.Range("D25").Value = "'4.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "'1.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.79%  "

# Row 27
$ws.Range("D27").Value = "'9.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.39%  "

# Row 28
$ws.Range("E28").Value = "  +1.46%  "

# Row 29
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("D30").Value = "'522.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.89%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0898"
$ws.Range("E31").Value = "  -1.38%  "

# Row 32
$ws.Range("D32").Value = "'7.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.41%  "

# Row 33
$ws.Range("E33").Value = "  +0.14%  "

# Row 34
$ws.Range("E34").Value = "  +0.49%  "

# Row 35
$ws.Range("E35").Value = "  -0.08%  "

# Row 36
$ws.Range("D36").Value = "'163.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.37%  "

# Row 37
$ws.Range("D37").Value = "'0.120"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.44%  "

# Row 38
$ws.Range("D38").Value = "'18.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.12%  "

# Row 39
$ws.Range("D39").Value = "'18.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "

# Row 40
$ws.Range("D40").Value = "'1.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.01%  "

# Row 41
$ws.Range("E41").Value = "  -0.85%  "

# Row 42
$ws.Range("E42").Value = "  +0.03%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "'4.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "

# Row 44
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.327"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.66%  "

# Row 45
$ws.Range("D45").Value = "'2.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.56%  "

# Row 46
$ws.Range("D46").Value = "'151.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.59%  "

# Row 47
$ws.Range("D47").Value = "'3.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.71%  "

# Row 48
$ws.Range("D48").Value = "'0.520"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "

# Row 49
$ws.Range("B49").Value = "Optimism"
$ws.Range("C49").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D49").Value = "'1.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0742"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.46%  "

# Row 51
$ws.Range("D51").Value = "'0.582"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.83%  "
